# Add IP on the file
# Populate the new IP-address / network-name columns (E:F) for rows 6-9
# and tidy up the view the way the author left it after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ipData = @(
    @{ Row = 6; Ip = '"http://192.168.1.18:8080"';  Name = "WiFi Anna" },
    @{ Row = 7; Ip = '"http://172.17.0.1:8080"';     Name = "Docker Container" },
    @{ Row = 8; Ip = '"http://172.20.10.3:8080"';    Name = "Iphone Anna hotspot" },
    @{ Row = 9; Ip = 'http://172.20.228.53:8080"';   Name = "Ubuntu macchina Virtuale" }
)

foreach ($entry in $ipData) {
    $eCell = $ws.Cells.Item($entry.Row, 5)
    $eCell.Value = $entry.Ip
    $eCell.Font.Name = "Consolas"
    $eCell.Font.Size = 8
    $eCell.Font.Color = 7901646
    $eCell.VerticalAlignment = -4108

    $fCell = $ws.Cells.Item($entry.Row, 6)
    $fCell.Value = $entry.Name
    $fCell.VerticalAlignment = -4108
}

# Autofit the new column so the IP/network labels are fully visible
$ws.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

# Leave the view where the author ended up after typing the data in F10
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("F10").Select()
